# Add scheduler in coding sections
# Reorders a handful of Item Name / UOM rows (within the same BRAND group)
# and fills in Total Ordered / Estimated Sales for row 2.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Item Name (column D) reorders ---
$ws.Range("D4").Value2  = "Dinafex 60mg Tablet"
$ws.Range("D5").Value2  = "Dinafex 180mg Tablet"
$ws.Range("D6").Value2  = "Dinafex 120mg Tablet"

$ws.Range("D8").Value2  = "Etorix 120mg Tablet"
$ws.Range("D9").Value2  = "Etorix 90mg Tablet"

$ws.Range("D12").Value2 = "Flucloxin 500mg Capsule - 36's"
$ws.Range("D13").Value2 = "Flucloxin 500mg Capsule"

$ws.Range("D15").Value2 = "Ketonic 10mg Tablet"
$ws.Range("D16").Value2 = "Ketonic 30mg Injection"

$ws.Range("D18").Value2 = "Kynol D 25mg Tablet"
$ws.Range("D19").Value2 = "Kynol TR 200mg Capsule"

$ws.Range("D26").Value2 = "Zithrox 30ml Dry Suspension"
$ws.Range("D27").Value2 = "Zithrox 500mg Tablet"

# --- UOM (column E) reorders ---
$ws.Range("E8").Value2  = "20's"
$ws.Range("E9").Value2  = "30's"

$ws.Range("E12").Value2 = "36 's"
$ws.Range("E13").Value2 = "30 's"

$ws.Range("E15").Value2 = "20's"
$ws.Range("E16").Value2 = "5 's"

$ws.Range("E18").Value2 = "60 's"
$ws.Range("E19").Value2 = "30 's"

$ws.Range("E21").Value2 = "30 's"

$ws.Range("E26").Value2 = "30ml"
$ws.Range("E27").Value2 = "6 's"

# --- Totals for row 2 (Biltin 20mg Tablet 30's) ---
$ws.Range("F2").Value2 = 9470
$ws.Range("G2").Value2 = 3194515
